# Update June 2022 figures on the "Table 1" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revise a few previously-reported monthly totals for 2022 (row 2).
$ws.Range("C2").Value2 = 13680   # NOV
$ws.Range("D2").Value2 = 15336   # DEC
$ws.Range("G2").Value2 = 16087   # MAR

# Report newly available monthly totals for 2022 (APR, MAY, JUN).
$ws.Range("H2").Value2 = 12635   # APR
$ws.Range("I2").Value2 = 19080   # MAY
$ws.Range("J2").Value2 = 11258   # JUN

# Give the newly-populated (and still-empty trailing) cells the same
# number format / alignment as the rest of the 2022 row.
$ws.Range("G2").Copy()
$ws.Range("H2:M2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Leave the active cell where the author left off editing.
[void]$ws.Range("J3").Select()
